# Update simulation results on Sheet1 (rows 2-25, columns B,D,E,F,G,I,K,L,M,O)
# "case with 380 kV done" -- refreshed pl_mw.xlsx result values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row data keyed by worksheet row number; each entry lists the new value
# for every column that changed in that row.
$rowUpdates = @(
    ,@{ Row = 2; Values = @(0.5563414777818139, 0.03501980377711078, 0.2865979018945954, 0.792792929073606, 0.002439753835041164, 1.133782779201486, 0.4491800626102247, 0.1500470941605272, 0.1504785850134134, 2.764484150819413) }
    ,@{ Row = 3; Values = @(0.535161830774598, 0.03300260780376618, 0.2884919792734988, 0.7921491412657744, 0.002442200790296443, 1.146193432821434, 0.3954152233587536, 0.1436718949989881, 0.1445082727467764, 2.775940267536797) }
    ,@{ Row = 4; Values = @(0.5223439493359763, 0.03175285957945562, 0.2897245121466823, 0.7922498292998625, 0.002443784469829714, 1.154312427702539, 0.3622321539591553, 0.139824560359159, 0.1408993370808957, 2.784612015987179) }
    ,@{ Row = 5; Values = @(0.5171679435240719, 0.03124079020687986, 0.2902443098618696, 0.7924155900241985, 0.002444450320342284, 1.157746457610967, 0.3486674444699531, 0.1382736963325115, 0.1394430843685335, 2.78855757602733) }
    ,@{ Row = 6; Values = @(0.5163113479648302, 0.03115559389573974, 0.2903316818182607, 0.7924506481828075, 0.002444562123301687, 1.158324255291896, 0.3464124988782658, 0.1380172033717244, 0.1392021490525366, 2.789237600702776) }
    ,@{ Row = 7; Values = @(0.5222739514055377, 0.03174596488363335, 0.2897314512877944, 0.7922515597666049, 0.002443793366601832, 1.154358232132946, 0.3620493860160252, 0.1398035760870826, 0.1408796390107732, 2.784663560087708) }
    ,@{ Row = 8; Values = @(0.5490003095907809, 0.03432660975087742, 0.2872365656025514, 0.7924680160384625, 0.00244058072025022, 1.137958474223471, 0.4306779029202801, 0.1478350596409257, 0.1484082897532986, 2.768094293404602) }
    ,@{ Row = 9; Values = @(0.6028721424637808, 0.03929764919282519, 0.2828942761948852, 0.7968285442411371, 0.00243492258571777, 1.109754066029758, 0.5638740676170926, 0.1641137639364132, 0.1636187032716485, 2.748600619860497) }
    ,@{ Row = 10; Values = @(0.6433225920770553, 0.04289439795742567, 0.280036922338903, 0.8024345519511229, 0.002431153009239722, 1.091440067153947, 0.6608644744433718, 0.1763933848314849, 0.1750614894918527, 2.742212898196613) }
    ,@{ Row = 11; Values = @(0.6619096482742464, 0.04451844501387114, 0.2788088037199208, 0.8055072188069374, 0.002429521447526395, 1.08363061450644, 0.7047945479019404, 0.1820485959262896, 0.1803243285042342, 2.741032188819844) }
    ,@{ Row = 12; Values = @(0.6689743910215782, 0.04513166366063359, 0.2783540188130418, 0.8067459080411936, 0.002428915525297116, 1.080748351460088, 0.7214016037051749, 0.1841999463289739, 0.1823253870851786, 2.74083328348874) }
    ,@{ Row = 13; Values = @(0.6674517138866918, 0.04499967534352578, 0.2784515083525316, 0.8064757923178689, 0.002429045492513169, 1.081365762560022, 0.7178262487886684, 0.1837361782593661, 0.1818940634968769, 2.74086508023683) }
    ,@{ Row = 14; Values = @(0.6624903461832048, 0.04456893057049882, 0.278771182521671, 0.8056076206237179, 0.002429471359445672, 1.083391985788218, 0.7061613922821834, 0.182225392083879, 0.1804887944725735, 2.741010849613673) }
    ,@{ Row = 15; Values = @(0.6594547657344947, 0.04430485526280847, 0.2789683295314023, 0.8050856256294878, 0.002429733765322268, 1.08464287389647, 0.6990126205948854, 0.181301271058274, 0.1796290823636681, 2.741132465023838) }
    ,@{ Row = 16; Values = @(0.6421115822341221, 0.0427880161281351, 0.2801186227512955, 0.8022442599287558, 0.002431261305765597, 1.091960929924955, 0.6579896224750428, 0.176025186501775, 0.1747186959520519, 2.742324785457072) }
    ,@{ Row = 17; Values = @(0.6315193618924297, 0.04185435710462571, 0.2808426315406727, 0.8006349956070835, 0.002432219682641048, 1.096583927595461, 0.6327737670859221, 0.1728061200915079, 0.171720959115234, 2.743498172454906) }
    ,@{ Row = 18; Values = @(0.6254445292510411, 0.04131620088870847, 0.2812658135881969, 0.7997585659597704, 0.002432778753688914, 1.099292056863106, 0.6182523106011502, 0.1709611128807325, 0.1700021559754887, 2.744335439927653) }
    ,@{ Row = 19; Values = @(0.6233907220727417, 0.04113379546210183, 0.2814102563555416, 0.7994702678365897, 0.002432969393507226, 1.100217416727553, 0.6133325339497446, 0.1703375470516448, 0.1694211327992221, 2.744646806669181) }
    ,@{ Row = 20; Values = @(0.6326451110021196, 0.04195386491626607, 0.2807648610801075, 0.8008012148887786, 0.002432116851108613, 1.096086719664093, 0.635459902842598, 0.1731481216815496, 0.1720395137228508, 2.743356458098958) }
    ,@{ Row = 21; Values = @(0.6639469123864785, 0.04469549908745307, 0.2786770078273728, 0.8058605844502438, 0.002429345949080446, 1.082794799309966, 0.7095884180265841, 0.1826688802128018, 0.1809013362589837, 2.740961296361576) }
    ,@{ Row = 22; Values = @(0.6845570896919355, 0.04647696965301407, 0.2773723612942078, 0.809605132978902, 0.002427604428546829, 1.07454495605484, 0.7578703113345, 0.1889485826800268, 0.1867404117456886, 2.740842702202173) }
    ,@{ Row = 23; Values = @(0.673543251902629, 0.0455271212264492, 0.2780632073641245, 0.8075665218257058, 0.002428527576460365, 1.078908052324515, 0.7321167497869965, 0.1855917737444059, 0.183619695463257, 2.740773572982363) }
    ,@{ Row = 24; Values = @(0.6321361137667907, 0.04190888174082374, 0.2807999994867745, 0.8007259152746116, 0.002432163316070953, 1.096311350668223, 0.6342455772024778, 0.1729934850727659, 0.1718954807358344, 2.743420020455034) }
    ,@{ Row = 25; Values = @(0.5881440087123053, 0.03796253167770658, 0.2840103373643486, 0.7952272294687575, 0.002436384944236559, 1.116961016995891, 0.5638740676170926, 0.1596536033947163, 0.1594565376865802, 2.75248154467522) }
)

$columns = @("B", "D", "E", "F", "G", "I", "K", "L", "M", "O")

foreach ($update in $rowUpdates) {
    $r = $update.Row
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $col = $columns[$i]
        $ws.Range("$col$r").Value = [double]$update.Values[$i]
    }
}